$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$wsMCF = $wb.Worksheets.Item("MCF")

# Update the "last updated" date on the About sheet (C1)
$wsAbout.Range("C1").Value = (Get-Date -Year 2024 -Month 4 -Day 10 -Hour 0 -Minute 0 -Second 0)

# Update Maximum Capacity Factor values on MCF sheet to 1
$wsMCF.Range("B2").Value = 1
$wsMCF.Range("B3").Value = 1
$wsMCF.Range("B4").Value = 1
$wsMCF.Range("B6").Value = 1
$wsMCF.Range("B10").Value = 1
$wsMCF.Range("B11").Value = 1
$wsMCF.Range("B12").Value = 1
$wsMCF.Range("B13").Value = 1
$wsMCF.Range("B14").Value = 1
$wsMCF.Range("B16").Value = 1
$wsMCF.Range("B17").Value = 1
$wsMCF.Range("B18").Value = 1

$wsMCF.Activate()
$wsMCF.Range("B17").Select()
